$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Preços" -> "Preços com Fator 3"
$ws.Range("C1").Value = "Preços com Fator 3"

# Part numbers (keep them as text, not auto-converted numbers).
# Leading apostrophe forces Excel to store the value as text; resetting the
# style back to "Normal" afterwards drops the text-prefix formatting so the
# cell itself stays a plain shared-string cell like the original.
$ws.Range("A2").Value = "'111954114"
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").Value = "'111937242"
$ws.Range("A3").Style = "Normal"

# Quantities / prices
$ws.Range("C2").Value = 269.7
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 4794

# Total label
$ws.Range("C4").Value = "Total: 5063.7"
